$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (changed) date in column C for all data rows (2-301)
#    from serial 45182 (2023-09-13) to serial 45184 (2023-09-15).
$ws.Range("C2:C301").Value = 45184

# 2. Row 301 gains an explicit row height (matches the other data rows).
$ws.Rows.Item(301).RowHeight = 15

# 3. Append new row 302: A 42721-2023
$ws.Rows.Item(302).RowHeight = 15
$ws.Range("A302").Value = "A 42721-2023"
$ws.Range("B302").NumberFormat = "YYYY-MM-DD"
$ws.Range("B302").Value = 45181
$ws.Range("C302").NumberFormat = "YYYY-MM-DD"
$ws.Range("C302").Value = 45184
$ws.Range("D302").Value = "NORRBOTTENS LÄN"
$ws.Range("E302").Value = "ÖVERKALIX"
$ws.Range("G302").Value = 13.1
$ws.Range("H302").Value = 0
$ws.Range("I302").Value = 0
$ws.Range("J302").Value = 0
$ws.Range("K302").Value = 0
$ws.Range("L302").Value = 0
$ws.Range("M302").Value = 0
$ws.Range("N302").Value = 0
$ws.Range("O302").Value = 0
$ws.Range("P302").Value = 0
$ws.Range("Q302").Value = 0
$ws.Range("R302").WrapText = $true

# 4. Append new row 303: A 43099-2023
$ws.Range("A303").Value = "A 43099-2023"
$ws.Range("B303").NumberFormat = "YYYY-MM-DD"
$ws.Range("B303").Value = 45182
$ws.Range("C303").NumberFormat = "YYYY-MM-DD"
$ws.Range("C303").Value = 45184
$ws.Range("D303").Value = "NORRBOTTENS LÄN"
$ws.Range("E303").Value = "ÖVERKALIX"
$ws.Range("G303").Value = 6.1
$ws.Range("H303").Value = 0
$ws.Range("I303").Value = 0
$ws.Range("J303").Value = 0
$ws.Range("K303").Value = 0
$ws.Range("L303").Value = 0
$ws.Range("M303").Value = 0
$ws.Range("N303").Value = 0
$ws.Range("O303").Value = 0
$ws.Range("P303").Value = 0
$ws.Range("Q303").Value = 0
$ws.Range("R303").WrapText = $true
